# Add two new columns (I = "I0", J = "IF") to the sheet, mirroring the
# existing header style used in row 1, and fill in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers - reuse the same formatting already used by the other header
# cells (row 1, style index 1: bold, centered, thin border) by copying it
# over rather than constructing a brand new style entry.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Row -> value map (I and J share the same value on every data row)
$values = @{
    2=9; 3=9; 4=9; 5=9; 6=9; 7=9; 8=9; 9=9; 10=9;
    11=9; 12=9; 13=9; 14=9; 15=9; 16=9; 17=10; 18=9; 19=9; 20=9;
    21=9; 22=9; 23=9; 24=9; 25=9; 26=9; 27=9; 28=9; 29=9; 30=9;
    31=9; 32=9; 33=9; 34=9; 35=9; 36=9; 37=9; 38=9; 39=9; 40=9;
    41=7; 42=7; 43=9; 44=9; 45=9; 46=10; 47=9; 48=9; 49=9; 50=9;
    51=9; 52=7; 53=7; 54=9; 55=9; 56=9; 57=9; 58=9; 59=9; 60=9;
    61=9; 62=9; 63=9; 64=9; 65=9; 66=9; 67=9; 68=9; 69=10; 70=9;
    71=9; 72=9; 73=9; 74=9; 75=9; 76=9; 77=9; 78=9; 79=9; 80=9;
    81=7; 82=9; 83=9; 84=5; 85=8; 86=7; 87=4; 88=4; 89=3
}

foreach ($r in $values.Keys) {
    $v = $values[$r]
    $ws.Cells.Item($r, 9).Value = $v
    $ws.Cells.Item($r, 10).Value = $v
}
